$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.27%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.45%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.037"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.86%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07981"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.63%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.906"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.28%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.803"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.051"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.40%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.989"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.51%"

$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9224"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.67%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1300"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "29.80%"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1853"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.93%"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09889"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "8.63%"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03515"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.53%"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09846"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.58%"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001386"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.54%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005799"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.38%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.506"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.13%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.27%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.26%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.043"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.05%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04497"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.34%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.46%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004787"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.73%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001251"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.01%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003001"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-33.38%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01881"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.74%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04718"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.71%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007544"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.33%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01022"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "30.17%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.78%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002112"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.00%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01057"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.35%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006242"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.49%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "64.15%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-12.48%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
